$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (E1:J1) - new mobility category labels
$headers = @("Retail","Grocery","Parks","Transit","Workplace","Residential")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $ws.Cells.Item(1, 5 + $col).Value = $headers[$col]
}

# Mobility data for rows 2:35, columns E:J
$data = @(
    @(-2,8,-11,-28,-10,6),
    @(-34,-26,-43,-69,-37,15),
    @(-13,9,-9,-36,-15,9),
    @(-11,-6,-4,-47,-23,14),
    @(-1,8,-3,-30,-16,8),
    @(-6,1,-9,-37,-18,13),
    @(-14,-6,-16,-27,-16,9),
    @(-1,6,1,-29,-11,8),
    @(-10,-4,-14,-38,-21,14),
    @(-13,-1,-10,-34,-19,11),
    @(-7,1,2,-19,-7,7),
    @(-4,-2,21,0,-17,13),
    @(-26,-8,-71,-42,-31,14),
    @(-11,3,-12,-27,-14,8),
    @(-3,13,-8,-29,-20,9),
    @(-20,-11,-14,-33,-20,11),
    @(-12,9,-7,-29,-11,8),
    @(-21,9,-20,-9,-18,10),
    @(-27,-8,-20,-45,-22,13),
    @(-20,-11,-16,-36,-22,11),
    @(-20,-11,-10,-31,-14,9),
    @(-9,-3,-5,-41,-15,9),
    @(-20,-12,-25,-50,-15,11),
    @(-3,8,-7,-28,-13,9),
    @(-17,-7,-10,-37,-19,12),
    @(-15,-2,-7,-47,-25,12),
    @(-13,1,5,-41,-19,10),
    @(-23,-17,-32,-56,-29,13),
    @(-23,3,-2,-28,-23,13),
    @(-23,-12,-18,-36,-18,10),
    @(-13,-4,-17,-53,-15,12),
    @(-20,-7,-31,-47,-19,11),
    @(1,13,2,-29,-12,10),
    @(-7,4,-9,-31,-19,8),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($col = 0; $col -lt $rowVals.Length; $col++) {
        $ws.Cells.Item($rowNum, 5 + $col).Value = $rowVals[$col]
    }
}

# Column widths (bestFit-style autosizing on the numeric mobility columns)
# (ColumnWidth is "characters"; the engine adds ~0.8333 padding when it
# stores the XML <col width>, so back the input off by that much to land
# on the saved widths from the diff: ~10.5 / ~10.5 / 11.)
$ws.Columns.Item(4).ColumnWidth = 9.666666666666666
$ws.Columns.Item(9).ColumnWidth = 9.666666666666666
$ws.Columns.Item(10).ColumnWidth = 10.166666666666666

# Selection left on H1, matching the saved view state
$ws.Range("H1").Select()
